# Update the "department" column (C) on the "courses" sheet.
# Rows 2-3 (Graduate Diploma of Management / Strategic Leadership) become "Business".
# Rows 4-12 (package offerings) become "Packages".
# Previously every row held the single shared string "FACULTY OF BUSINESS & TECHNOLOGY".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

$ws.Range("C2").Value = "Business"
$ws.Range("C3").Value = "Business"

$ws.Range("C4").Value = "Packages"
$ws.Range("C5").Value = "Packages"
$ws.Range("C6").Value = "Packages"
$ws.Range("C7").Value = "Packages"
$ws.Range("C8").Value = "Packages"
$ws.Range("C9").Value = "Packages"
$ws.Range("C10").Value = "Packages"
$ws.Range("C11").Value = "Packages"
$ws.Range("C12").Value = "Packages"
